$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (project "Acacia Breeze"): opening/closing dates shift ---
$ws.Cells.Item(2, 9).Value  = 45762   # I2 Application opening date
$ws.Cells.Item(2, 10).Value = 45797   # J2 Application closing date

# --- New project rows (3-7): fill column-by-column so new shared strings are
#     interned in the same order the source workbook used. ---

# Column A - Project Name
$ws.Cells.Item(3, 1).Value = "Serenity Grove"
$ws.Cells.Item(5, 1).Value = "Verdant Rise"
$ws.Cells.Item(6, 1).Value = "Coastal Haven"
$ws.Cells.Item(7, 1).Value = "Skyline Nest"
$ws.Cells.Item(4, 1).Value = "Amber Heights"

# Column B - Neighborhood
$ws.Cells.Item(4, 2).Value = "Tampines"
$ws.Cells.Item(5, 2).Value = "Bukit Panjang"
$ws.Cells.Item(6, 2).Value = "Pasir Ris"
$ws.Cells.Item(7, 2).Value = "Punggol"
$ws.Cells.Item(3, 2).Value = "Yishun"

# Column K - Manager
$ws.Cells.Item(3, 11).Value = "Michael"
$ws.Cells.Item(4, 11).Value = "Jessica"
$ws.Cells.Item(5, 11).Value = "James"
$ws.Cells.Item(6, 11).Value = "Frank"
$ws.Cells.Item(7, 11).Value = "Kelly"

# Remaining columns (no new shared strings introduced)
$data = @(
    @{ Row=3; C="2-Room"; D=3; E=300000; F="3-Room"; G=2; H=400000; I=45767; J=45798; L=4 },
    @{ Row=4; C="2-Room"; D=4; E=250000; F="3-Room"; G=4; H=325000; I=45698; J=45738; L=6 },
    @{ Row=5; C="2-Room"; D=5; E=375000; F="3-Room"; G=3; H=450000; I=45777; J=45800; L=5 },
    @{ Row=6; C="2-Room"; D=6; E=600000; F="3-Room"; G=4; H=700000; I=45748; J=45801; L=7 },
    @{ Row=7; C="2-Room"; D=7; E=500000; F="3-Room"; G=3; H=600000; I=45759; J=45802; L=5 }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 3).Value  = $row.C
    $ws.Cells.Item($r, 4).Value  = $row.D
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 6).Value  = $row.F
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 12).Value = $row.L
}

# Apply the same date formatting used by row 2 (I/J columns) to the new rows
$ws.Range("I2:J2").Copy() | Out-Null
$ws.Range("I3:J7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Column width adjustments (closing in on the author's layout) ---
$ws.Columns.Item(2).ColumnWidth = 15.3
$ws.Columns.Item(9).ColumnWidth = 13.3
$ws.Columns.Item(10).ColumnWidth = 14.8

# --- Update selection to match the author's last active cell ---
$ws.Range("H5").Select() | Out-Null
